$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.029686947197129
$ws.Cells.Item(2, 4).Value = 1.032986501582317
$ws.Cells.Item(2, 5).Value = 1.042840500166636
$ws.Cells.Item(2, 6).Value = 1.049611902960549
$ws.Cells.Item(2, 9).Value = 1.031518270347117
$ws.Cells.Item(2, 10).Value = 1.03483234535068
$ws.Cells.Item(2, 11).Value = 1.035790151433865
$ws.Cells.Item(2, 12).Value = 1.04561605158057
$ws.Cells.Item(2, 13).Value = 1.052368477587838
$ws.Cells.Item(2, 14).Value = 1.015600469996057
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.030702254613439
$ws.Cells.Item(3, 4).Value = 1.033742472218859
$ws.Cells.Item(3, 5).Value = 1.043877789915342
$ws.Cells.Item(3, 6).Value = 1.050793132461464
$ws.Cells.Item(3, 9).Value = 1.031699517605085
$ws.Cells.Item(3, 10).Value = 1.035488453425146
$ws.Cells.Item(3, 11).Value = 1.036355149853742
$ws.Cells.Item(3, 12).Value = 1.046463662223838
$ws.Cells.Item(3, 13).Value = 1.053361035782542
$ws.Cells.Item(3, 14).Value = 1.015818769069912
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.031359115764914
$ws.Cells.Item(4, 4).Value = 1.034230949566993
$ws.Cells.Item(4, 5).Value = 1.044549290760344
$ws.Cells.Item(4, 6).Value = 1.051558062660346
$ws.Cells.Item(4, 9).Value = 1.031814712724328
$ws.Cells.Item(4, 10).Value = 1.035912317439308
$ws.Cells.Item(4, 11).Value = 1.036719396259619
$ws.Cells.Item(4, 12).Value = 1.047011825715351
$ws.Cells.Item(4, 13).Value = 1.054003309014135
$ws.Cells.Item(4, 14).Value = 1.015959752729004
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.031635233739195
$ws.Cells.Item(5, 4).Value = 1.034436140355133
$ws.Cells.Item(5, 5).Value = 1.044831662514357
$ws.Cells.Item(5, 6).Value = 1.05187978198264
$ws.Cells.Item(5, 9).Value = 1.031862641353756
$ws.Cells.Item(5, 10).Value = 1.036090346350398
$ws.Cells.Item(5, 11).Value = 1.036872202790693
$ws.Cells.Item(5, 12).Value = 1.04724220226197
$ws.Cells.Item(5, 13).Value = 1.054273325798512
$ws.Cells.Item(5, 14).Value = 1.016018957372081
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.031681593611703
$ws.Cells.Item(6, 4).Value = 1.034470583086679
$ws.Cells.Item(6, 5).Value = 1.044879078236756
$ws.Cells.Item(6, 6).Value = 1.051933808439908
$ws.Cells.Item(6, 9).Value = 1.031870659482139
$ws.Cells.Item(6, 10).Value = 1.036120228590179
$ws.Cells.Item(6, 11).Value = 1.036897840741838
$ws.Cells.Item(6, 12).Value = 1.047280879305379
$ws.Cells.Item(6, 13).Value = 1.054318663102477
$ws.Cells.Item(6, 14).Value = 1.016028894282937
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.031362805370714
$ws.Cells.Item(7, 4).Value = 1.034233691986078
$ws.Cells.Item(7, 5).Value = 1.044553063539785
$ws.Cells.Item(7, 6).Value = 1.051562360929042
$ws.Cells.Item(7, 9).Value = 1.031815355111951
$ws.Cells.Item(7, 10).Value = 1.035914696913249
$ws.Cells.Item(7, 11).Value = 1.036721439336543
$ws.Cells.Item(7, 12).Value = 1.047014904298825
$ws.Cells.Item(7, 13).Value = 1.054006916971626
$ws.Cells.Item(7, 14).Value = 1.01596054407949
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.030030098022241
$ws.Cells.Item(8, 4).Value = 1.033242126997777
$ws.Cells.Item(8, 5).Value = 1.043190993806177
$ws.Cells.Item(8, 6).Value = 1.050010981682896
$ws.Cells.Item(8, 9).Value = 1.031579954915187
$ws.Cells.Item(8, 10).Value = 1.03505422113547
$ws.Cells.Item(8, 11).Value = 1.035981373514926
$ws.Cells.Item(8, 12).Value = 1.045902566989402
$ws.Cells.Item(8, 13).Value = 1.052703912735822
$ws.Cells.Item(8, 14).Value = 1.015674301100698
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.027680847730873
$ws.Cells.Item(9, 4).Value = 1.031489642223545
$ws.Cells.Item(9, 5).Value = 1.040793198012491
$ws.Cells.Item(9, 6).Value = 1.047281812775156
$ws.Cells.Item(9, 9).Value = 1.031149210022305
$ws.Cells.Item(9, 10).Value = 1.033532743815243
$ws.Cells.Item(9, 11).Value = 1.034667003911144
$ws.Cells.Item(9, 12).Value = 1.043940218635538
$ws.Cells.Item(9, 13).Value = 1.050408021768254
$ws.Cells.Item(9, 14).Value = 1.015167839278109
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.026114102812421
$ws.Cells.Item(10, 4).Value = 1.030317853511226
$ws.Cells.Item(10, 5).Value = 1.039196253063211
$ws.Cells.Item(10, 6).Value = 1.045465424471224
$ws.Cells.Item(10, 9).Value = 1.030851356800987
$ws.Cells.Item(10, 10).Value = 1.032514937785506
$ws.Cells.Item(10, 11).Value = 1.033783877330904
$ws.Cells.Item(10, 12).Value = 1.04263047053748
$ws.Cells.Item(10, 13).Value = 1.048877541620031
$ws.Cells.Item(10, 14).Value = 1.014828817033623
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.025435545240738
$ws.Cells.Item(11, 4).Value = 1.029809642644968
$ws.Cells.Item(11, 5).Value = 1.038505135150291
$ws.Cells.Item(11, 6).Value = 1.044679630587578
$ws.Cells.Item(11, 9).Value = 1.030719852761021
$ws.Cells.Item(11, 10).Value = 1.032073391519784
$ws.Cells.Item(11, 11).Value = 1.033399847618626
$ws.Cells.Item(11, 12).Value = 1.042062976529541
$ws.Cells.Item(11, 13).Value = 1.048214852759457
$ws.Cells.Item(11, 14).Value = 1.014681690407014
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.025183476358597
$ws.Cells.Item(12, 4).Value = 1.029620748106988
$ws.Cells.Item(12, 5).Value = 1.038248478686742
$ws.Cells.Item(12, 6).Value = 1.04438785868248
$ws.Cells.Item(12, 9).Value = 1.030670626304502
$ws.Cells.Item(12, 10).Value = 1.031909257031311
$ws.Cells.Item(12, 11).Value = 1.033256957347855
$ws.Cells.Item(12, 12).Value = 1.041852129195152
$ws.Cells.Item(12, 13).Value = 1.047968703273455
$ws.Cells.Item(12, 14).Value = 1.014626991858783
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.025237547008217
$ws.Cells.Item(13, 4).Value = 1.029661272155768
$ws.Cells.Item(13, 5).Value = 1.038303529858806
$ws.Cells.Item(13, 6).Value = 1.044450439913797
$ws.Cells.Item(13, 9).Value = 1.030681202736859
$ws.Cells.Item(13, 10).Value = 1.031944470070931
$ws.Cells.Item(13, 11).Value = 1.033287618856049
$ws.Cells.Item(13, 12).Value = 1.041897359144268
$ws.Cells.Item(13, 13).Value = 1.048021503056738
$ws.Cells.Item(13, 14).Value = 1.014638727106867
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025414709607037
$ws.Cells.Item(14, 4).Value = 1.029794031056938
$ws.Cells.Item(14, 5).Value = 1.038483918715077
$ws.Cells.Item(14, 6).Value = 1.044655510451714
$ws.Cells.Item(14, 9).Value = 1.030715791436686
$ws.Cells.Item(14, 10).Value = 1.032059826651996
$ws.Cells.Item(14, 11).Value = 1.033388041250675
$ws.Cells.Item(14, 12).Value = 1.042045548941096
$ws.Cells.Item(14, 13).Value = 1.048194505913079
$ws.Cells.Item(14, 14).Value = 1.014677170010554
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.025523862354218
$ws.Cells.Item(15, 4).Value = 1.029875811996725
$ws.Cells.Item(15, 5).Value = 1.038595069595738
$ws.Cells.Item(15, 6).Value = 1.044781875336376
$ws.Cells.Item(15, 9).Value = 1.0307370523327
$ws.Cells.Item(15, 10).Value = 1.032130885147736
$ws.Cells.Item(15, 11).Value = 1.033449882409593
$ws.Cells.Item(15, 12).Value = 1.042136846304042
$ws.Cells.Item(15, 13).Value = 1.048301099029119
$ws.Cells.Item(15, 14).Value = 1.014700849439181
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.026159133323475
$ws.Cells.Item(16, 4).Value = 1.030351564596249
$ws.Cells.Item(16, 5).Value = 1.039242128069376
$ws.Cells.Item(16, 6).Value = 1.045517590039553
$ws.Cells.Item(16, 9).Value = 1.030860030982355
$ws.Cells.Item(16, 10).Value = 1.032544224286783
$ws.Cells.Item(16, 11).Value = 1.033809329805581
$ws.Cells.Item(16, 12).Value = 1.042668125564248
$ws.Cells.Item(16, 13).Value = 1.048921522482442
$ws.Cells.Item(16, 14).Value = 1.014838574444134
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.026557582502362
$ws.Cells.Item(17, 4).Value = 1.030649772887876
$ws.Cells.Item(17, 5).Value = 1.039648109720668
$ws.Cells.Item(17, 6).Value = 1.045979275393948
$ws.Cells.Item(17, 9).Value = 1.030936494699099
$ws.Cells.Item(17, 10).Value = 1.032803278978284
$ws.Cells.Item(17, 11).Value = 1.034034365503079
$ws.Cells.Item(17, 12).Value = 1.04300128541156
$ws.Cells.Item(17, 13).Value = 1.049310702699902
$ws.Cells.Item(17, 14).Value = 1.014924878003758
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.026789977038121
$ws.Cells.Item(18, 4).Value = 1.030823633714791
$ws.Cells.Item(18, 5).Value = 1.039884947632442
$ws.Cells.Item(18, 6).Value = 1.046248637589146
$ws.Cells.Item(18, 9).Value = 1.030980850402949
$ws.Cells.Item(18, 10).Value = 1.03295430115332
$ws.Cells.Item(18, 11).Value = 1.034165467715646
$ws.Cells.Item(18, 12).Value = 1.043195576753813
$ws.Cells.Item(18, 13).Value = 1.04953770679795
$ws.Cells.Item(18, 14).Value = 1.014975185781284
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.026869215183237
$ws.Cells.Item(19, 4).Value = 1.030882902341656
$ws.Cells.Item(19, 5).Value = 1.039965709294004
$ws.Cells.Item(19, 6).Value = 1.046340494841382
$ws.Cells.Item(19, 9).Value = 1.030995933105454
$ws.Cells.Item(19, 10).Value = 1.033005782232828
$ws.Cells.Item(19, 11).Value = 1.03421014348971
$ws.Cells.Item(19, 12).Value = 1.04326181911794
$ws.Cells.Item(19, 13).Value = 1.049615109651239
$ws.Cells.Item(19, 14).Value = 1.01499233406631
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.026514834140148
$ws.Cells.Item(20, 4).Value = 1.030617786129836
$ws.Cells.Item(20, 5).Value = 1.039604548041019
$ws.Cells.Item(20, 6).Value = 1.045929733796611
$ws.Cells.Item(20, 9).Value = 1.030928316131254
$ws.Cells.Item(20, 10).Value = 1.032775493140128
$ws.Cells.Item(20, 11).Value = 1.034010237570242
$ws.Cells.Item(20, 12).Value = 1.042965544149853
$ws.Cells.Item(20, 13).Value = 1.049268947154558
$ws.Cells.Item(20, 14).Value = 1.014915621719511
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.025362540258172
$ws.Cells.Item(21, 4).Value = 1.02975494024343
$ws.Cells.Item(21, 5).Value = 1.03843079716134
$ws.Cells.Item(21, 6).Value = 1.044595119340768
$ws.Cells.Item(21, 9).Value = 1.030705616413812
$ws.Cells.Item(21, 10).Value = 1.032025860451486
$ws.Cells.Item(21, 11).Value = 1.033358476105955
$ws.Cells.Item(21, 12).Value = 1.042001912252338
$ws.Cells.Item(21, 13).Value = 1.04814356081209
$ws.Cells.Item(21, 14).Value = 1.014665850890287
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.024637917658455
$ws.Cells.Item(22, 4).Value = 1.029211726939931
$ws.Cells.Item(22, 5).Value = 1.037693134479832
$ws.Cells.Item(22, 6).Value = 1.043756612741451
$ws.Cells.Item(22, 9).Value = 1.03056339786177
$ws.Cells.Item(22, 10).Value = 1.031553816082289
$ws.Cells.Item(22, 11).Value = 1.032947273164381
$ws.Cells.Item(22, 12).Value = 1.041395721598093
$ws.Cells.Item(22, 13).Value = 1.047436001050998
$ws.Cells.Item(22, 14).Value = 1.014508525618929
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.025022066221062
$ws.Cells.Item(23, 4).Value = 1.029499761490343
$ws.Cells.Item(23, 5).Value = 1.038084153009295
$ws.Cells.Item(23, 6).Value = 1.044201062456588
$ws.Cells.Item(23, 9).Value = 1.030638998875836
$ws.Cells.Item(23, 10).Value = 1.031804124006916
$ws.Cells.Item(23, 11).Value = 1.033165393596585
$ws.Cells.Item(23, 12).Value = 1.041717104799173
$ws.Cells.Item(23, 13).Value = 1.047811090490071
$ws.Cells.Item(23, 14).Value = 1.014591953659368
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.026534150328385
$ws.Cells.Item(24, 4).Value = 1.030632239813064
$ws.Cells.Item(24, 5).Value = 1.039624231578539
$ws.Cells.Item(24, 6).Value = 1.045952119298708
$ws.Cells.Item(24, 9).Value = 1.0309320124288
$ws.Cells.Item(24, 10).Value = 1.032788048611571
$ws.Cells.Item(24, 11).Value = 1.034021140430479
$ws.Cells.Item(24, 12).Value = 1.042981694196038
$ws.Cells.Item(24, 13).Value = 1.049287814682203
$ws.Cells.Item(24, 14).Value = 1.014919804333627
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.02828828661135
$ws.Cells.Item(25, 4).Value = 1.031943315580443
$ws.Cells.Item(25, 5).Value = 1.041412806856019
$ws.Cells.Item(25, 6).Value = 1.047986828825878
$ws.Cells.Item(25, 9).Value = 1.031262454025071
$ws.Cells.Item(25, 10).Value = 1.033926697956703
$ws.Cells.Item(25, 11).Value = 1.035008014859684
$ws.Cells.Item(25, 12).Value = 1.044447800534072
$ws.Cells.Item(25, 13).Value = 1.051001544546856
$ws.Cells.Item(25, 14).Value = 1.015299015750895
